# Add Send Email Module
# ----------------------
# The Login / AddAlbum / ChangePassword sheets keep exactly the same data
# they had before (only their view/selection state changes), so we just
# touch their selection state to reproduce the saved cursor positions, and
# autofit the new-ish column on AddAlbum. The main change is a brand new
# "SendEmail" worksheet with its own data + hyperlinks.

$wb = $excel.ActiveWorkbook

# --- Login sheet -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:B2").Select() | Out-Null

# --- AddAlbum sheet ----------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns("B").AutoFit() | Out-Null
$ws2.Range("D10").Select() | Out-Null

# --- ChangePassword sheet ----------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E10").Select() | Out-Null

# --- New SendEmail sheet -------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "SendEmail"

$ws4.Range("A1").Value = "To"
$ws4.Range("B1").Value = "Bcc"
$ws4.Range("C1").Value = "CC"
$ws4.Range("D1").Value = "Subject"
$ws4.Range("E1").Value = "Body"

$ws4.Range("A2").Value = "shaheen.akhtar@aimbeyond.com"
$ws4.Range("B2").Value = "shaheen.akhtar@aimbeyond.com"
$ws4.Range("C2").Value = "shaheen.akhtar@aimbeyond.com"
$ws4.Range("D2").Value = "Test"
$ws4.Range("E2").Value = "This is a test email."

$ws4.Hyperlinks.Add($ws4.Range("A2"), "mailto:shaheen.akhtar@aimbeyond.com") | Out-Null
$ws4.Hyperlinks.Add($ws4.Range("B2"), "mailto:shaheen.akhtar@aimbeyond.com") | Out-Null
$ws4.Hyperlinks.Add($ws4.Range("C2"), "mailto:shaheen.akhtar@aimbeyond.com") | Out-Null

# Re-apply the Hyperlink cell style after adding the links so the same
# style slot used elsewhere in the workbook gets reused instead of a new
# (duplicate) style being created.
$ws4.Range("A2").Style = "Hyperlink"
$ws4.Range("B2").Style = "Hyperlink"
$ws4.Range("C2").Style = "Hyperlink"

$ws4.Columns("A:C").AutoFit() | Out-Null
$ws4.Columns("E").AutoFit() | Out-Null

$ws4.Range("A3:E3").Select() | Out-Null
